$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting from the
# existing header cell H1 so they pick up the same bold/bordered/centered
# style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-15: new I (col 9) and J (col 10) numeric values.
$data = @{
    2  = @(6, 8)
    3  = @(1, 7)
    4  = @(2, 7)
    5  = @(1, 6)
    6  = @(1, 7)
    7  = @(1, 5)
    8  = @(5, 8)
    9  = @(2, 6)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 3)
    13 = @(1, 3)
    14 = @(1, 2)
    15 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
